$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 4: turn the old "Step 3 / DropDown" step into "Step 5 / DropDown"
#     with the dropdown's xpath (Page Object) and the selected value (Data Set) ---
$ws.Range("B4").Value = "Step 5"
$ws.Range("C4").Value = "//select[@id='landingTerm']"
$ws.Range("E4").Value = "1 Year Open 4%"
# Column E defaults new cells to style index 5 (inherited col style) - the target
# keeps E4 on the workbook's default/no style, so reset it explicitly.
$ws.Range("E4").Style = "Normal"

# --- Row 5 (new): a follow-up NavigateURL step back to a plain URL ---
$ws.Range("A5").Value = $ws.Range("A4").Value2
$ws.Range("B5").Value = "Step 6"
$ws.Range("D5").Value = "NavigateURL"
$ws.Range("E5").Value = "https://www.google.ca"

# Give E5 the same "hyperlink-look" style the old E3 cell used, since the
# hyperlink itself is about to move from E3 down to E5.
$ws.Range("E3").Copy()
$ws.Range("E5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the hyperlink (same target URL) from E3 to E5.
$ws.Range("E3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E5"), "https://tools.td.com/mortgage-payment-calculator/")

# Re-apply the hyperlink-look style, since adding the hyperlink nudges the style.
$ws.Range("E3").Copy()
$ws.Range("E5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Sheet cosmetics: widen column C for the new xpath text, move the active
#     selection to the newly added dropdown-value cell ---
$ws.Range("C1").EntireColumn.ColumnWidth = 37.166666666666664
$ws.Range("E4").Select()
